$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Demo Fund 2 / TSTF2 Port Co 3, amount update ---
$ws.Cells.Item(2, 1).Value = "Demo Fund 2"
$ws.Cells.Item(2, 2).Value = "TSTF2 Port Co 3"
$ws.Cells.Item(2, 4).Value = 100000000

# --- Row 3: Demo Fund 2 / TSTF2 Port Co 3, amount update ---
$ws.Cells.Item(3, 1).Value = "Demo Fund 2"
$ws.Cells.Item(3, 2).Value = "TSTF2 Port Co 3"
$ws.Cells.Item(3, 4).Value = 40000000

# --- Row 4: Demo Fund 2 / TSTF2 Port Co 4, new investment date, amount update ---
$ws.Cells.Item(4, 1).Value = "Demo Fund 2"
$ws.Cells.Item(4, 2).Value = "TSTF2 Port Co 4"
$ws.Cells.Item(4, 3).Value = 45214
$ws.Cells.Item(4, 4).Value = 280000000

# --- Row 5: Demo Fund 2 / TSTF2 Port Co 3, new date, amount, instrument ---
$ws.Cells.Item(5, 1).Value = "Demo Fund 2"
$ws.Cells.Item(5, 2).Value = "TSTF2 Port Co 3"
$ws.Cells.Item(5, 3).Value = 45717
$ws.Cells.Item(5, 4).Value = 100000000
# Quantity column is no longer a formula for this row - it's a manual
# negative adjustment entry now.
$ws.Cells.Item(5, 5).Value = -500000
# Instrument switches from CCPS back to Equity
$ws.Cells.Item(5, 7).Value = "Equity"
# Trailing helper cells J5/K5 are no longer used on this row
$ws.Cells.Item(5, 10).Clear() | Out-Null
$ws.Cells.Item(5, 11).Clear() | Out-Null

# The trailing near-empty row (row 6) is removed entirely.
$ws.Rows(6).Delete() | Out-Null

# Keep the hidden AutoFilter defined name in sync with the shrunk data range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$ALZ`$4"
    }
}

# Match the last recorded selection from the source edit.
$ws.Range("G5").Select() | Out-Null
